$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1160
$ws.Range("I20").Value = 1150
$ws.Range("J20").Value = 1200
$ws.Range("K20").Value = 1150
$ws.Range("L20").Value = 1200
$ws.Range("M20").Value = -920
$ws.Range("N20").Value = -1660

$ws.Range("H35").Value = 1160
$ws.Range("I35").Value = 1150
$ws.Range("J35").Value = 1200
$ws.Range("K35").Value = 1150
$ws.Range("L35").Value = 1200
$ws.Range("M35").Value = -771
$ws.Range("N35").Value = -1958

$ws.Range("H107").Value = 490.66666
$ws.Range("I107").Value = 461.75
$ws.Range("K107").Value = 461.75
$ws.Range("M107").Value = 1458.25

$ws.Range("H115").Value = 1010.3333
$ws.Range("I115").Value = 1010.3333
$ws.Range("K115").Value = 3030.9999
$ws.Range("M115").Value = -1463.9999

$ws.Range("H132").Value = 13016.777
$ws.Range("I132").Value = 13642.117
$ws.Range("J132").Value = 2386
$ws.Range("K132").Value = 40926.351
$ws.Range("L132").Value = 7158
$ws.Range("M132").Value = -38396.351
$ws.Range("N132").Value = -12218

$ws.Range("H137").Value = 2617.182
$ws.Range("I137").Value = 1846.75
$ws.Range("J137").Value = 3057.4285
$ws.Range("K137").Value = 5540.25
$ws.Range("L137").Value = 9172.2855
$ws.Range("M137").Value = -2990.25
$ws.Range("N137").Value = -14272.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3848324.5
$ws.Range("I32").Value = 519.86365
$ws.Range("K32").Value = 519.86365
$ws.Range("M32").Value = -232.86365

$ws.Range("H69").Value = 249999.5
$ws.Range("J69").Value = 249999.5
$ws.Range("L69").Value = 249999.5
$ws.Range("N69").Value = -251497.5

$ws.Range("H72").Value = 249999.5
$ws.Range("J72").Value = 249999.5
$ws.Range("L72").Value = 749998.5
$ws.Range("N72").Value = -757486.5

$ws.Range("H74").Value = 5456.4287
$ws.Range("I74").Value = 5016.6665
$ws.Range("K74").Value = 5016.6665
$ws.Range("M74").Value = -4142.6665

$ws.Range("H77").Value = 5456.4287
$ws.Range("I77").Value = 5016.6665
$ws.Range("K77").Value = 25083.3325
$ws.Range("M77").Value = -20715.3325

$ws.Range("H102").Value = 20840836
$ws.Range("I102").Value = 62502504
$ws.Range("K102").Value = 62502504
$ws.Range("M102").Value = -62500882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 249999.5
$ws.Range("J70").Value = 249999.5
$ws.Range("L70").Value = 249999.5
$ws.Range("N70").Value = -250585.5

$ws.Range("H73").Value = 249999.5
$ws.Range("J73").Value = 249999.5
$ws.Range("L73").Value = 249999.5
$ws.Range("N73").Value = -252027.5

$ws.Range("H94").Value = 403.5
$ws.Range("I94").Value = 394.72726
$ws.Range("K94").Value = 394.72726
$ws.Range("M94").Value = 56.27274

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6401.2573
$ws.Range("I31").Value = 4027.111
$ws.Range("K31").Value = 4027.111
$ws.Range("M31").Value = -3732.111

$ws.Range("H34").Value = 6401.2573
$ws.Range("I34").Value = 4027.111
$ws.Range("K34").Value = 4027.111
$ws.Range("M34").Value = -3825.111

$ws.Range("H58").Value = 2673.375
$ws.Range("I58").Value = 1770.2858
$ws.Range("K58").Value = 1770.2858
$ws.Range("M58").Value = -1567.2858

$ws.Range("H86").Value = 6950
$ws.Range("I86").Value = 6950
$ws.Range("K86").Value = 6950
$ws.Range("M86").Value = -5827

$ws.Range("H89").Value = 6950
$ws.Range("I89").Value = 6950
$ws.Range("K89").Value = 34750
$ws.Range("M89").Value = -29134

$ws.Range("H122").Value = 615.6
$ws.Range("I122").Value = 687.75
$ws.Range("J122").Value = 327
$ws.Range("K122").Value = 2063.25
$ws.Range("L122").Value = 981
$ws.Range("M122").Value = 386.75
$ws.Range("N122").Value = -5881

$ws.Range("H132").Value = 1507.6296
$ws.Range("I132").Value = 1373.3462
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4120.0386
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1590.0386
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 2673.375
$ws.Range("I136").Value = 1770.2858
$ws.Range("K136").Value = 5310.857400000001
$ws.Range("M136").Value = -2760.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 24.666666
$ws.Range("I2").Value = 17.857143
$ws.Range("J2").Value = 29
$ws.Range("K2").Value = 107.142858
$ws.Range("L2").Value = 174
$ws.Range("M2").Value = 5.857141999999996
$ws.Range("N2").Value = -400

$ws.Range("H22").Value = 746.6667
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 895
$ws.Range("K22").Value = 1350
$ws.Range("L22").Value = 2685
$ws.Range("M22").Value = -1181
$ws.Range("N22").Value = -3023

$ws.Range("H27").Value = 746.6667
$ws.Range("I27").Value = 450
$ws.Range("J27").Value = 895
$ws.Range("K27").Value = 1350
$ws.Range("L27").Value = 2685
$ws.Range("M27").Value = -1248
$ws.Range("N27").Value = -2889

$ws.Range("H38").Value = 291.20834
$ws.Range("I38").Value = 303.65216
$ws.Range("J38").Value = 5
$ws.Range("K38").Value = 910.9564799999999
$ws.Range("L38").Value = 15
$ws.Range("M38").Value = -563.9564799999999
$ws.Range("N38").Value = -709

$ws.Range("H68").Value = 600.75
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 634.3333
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 1902.9999
$ws.Range("M68").Value = -689
$ws.Range("N68").Value = -3524.9999

$ws.Range("H71").Value = 600.75
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 634.3333
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 5708.9997
$ws.Range("M71").Value = -444
$ws.Range("N71").Value = -13820.9997

$ws.Range("H107").Value = 1500
$ws.Range("J107").Value = 1500
$ws.Range("L107").Value = 4500
$ws.Range("N107").Value = -8340

$ws.Range("H140").Value = 4078.6667
$ws.Range("I140").Value = 3572.25
$ws.Range("K140").Value = 10716.75
$ws.Range("M140").Value = -5536.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = 0

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = 0

$ws.Range("H80").Value = 4537.625
$ws.Range("I80").Value = 4450.1665
$ws.Range("K80").Value = 4450.1665
$ws.Range("M80").Value = -3452.1665

$ws.Range("H83").Value = 4537.625
$ws.Range("I83").Value = 4450.1665
$ws.Range("K83").Value = 22250.8325
$ws.Range("M83").Value = -17258.8325

$ws.Range("H97").Value = 1281.6666
$ws.Range("I97").Value = 1422
$ws.Range("J97").Value = 1106.25
$ws.Range("K97").Value = 1422
$ws.Range("L97").Value = 1106.25
$ws.Range("M97").Value = -926
$ws.Range("N97").Value = -2098.25

$ws.Range("H122").Value = 2938.6667
$ws.Range("I122").Value = 800
$ws.Range("K122").Value = 2400
$ws.Range("M122").Value = 50

$ws.Range("H132").Value = 51043.145
$ws.Range("I132").Value = 58561.5
$ws.Range("J132").Value = 5933
$ws.Range("K132").Value = 175684.5
$ws.Range("L132").Value = 17799
$ws.Range("M132").Value = -173154.5
$ws.Range("N132").Value = -22859

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 9000
$ws.Range("I17").Value = 9000
$ws.Range("K17").Value = 9000
$ws.Range("M17").Value = -8830

$ws.Range("H40").Value = 9377.25
$ws.Range("I40").Value = 9377.25
$ws.Range("K40").Value = 9377.25
$ws.Range("M40").Value = -9241.25

$ws.Range("H82").Value = 3374.5833
$ws.Range("J82").Value = 6600
$ws.Range("L82").Value = 6600
$ws.Range("N82").Value = -7322

$ws.Range("H85").Value = 3374.5833
$ws.Range("J85").Value = 6600
$ws.Range("L85").Value = 6600
$ws.Range("N85").Value = -9096

$ws.Range("H93").Value = 1513.7778
$ws.Range("I93").Value = 1260.5714
$ws.Range("K93").Value = 1260.5714
$ws.Range("M93").Value = -12.57140000000004

$ws.Range("H122").Value = 7502
$ws.Range("I122").Value = 7502
$ws.Range("K122").Value = 22506
$ws.Range("M122").Value = -20056

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 32714.5
$ws.Range("J46").Value = 32714.5
$ws.Range("L46").Value = 32714.5
$ws.Range("N46").Value = -33176.5

$ws.Range("H81").Value = 715.6
$ws.Range("I81").Value = 715.6
$ws.Range("K81").Value = 1431.2
$ws.Range("M81").Value = -370.2

$ws.Range("H84").Value = 715.6
$ws.Range("I84").Value = 715.6
$ws.Range("K84").Value = 7156
$ws.Range("M84").Value = -1852

$ws.Range("H107").Value = 27778386
$ws.Range("I107").Value = 41667120
$ws.Range("K107").Value = 125001360
$ws.Range("M107").Value = -124999440

$ws.Range("H113").Value = 1110.7142
$ws.Range("I113").Value = 986.2
$ws.Range("K113").Value = 2958.6
$ws.Range("M113").Value = -788.6000000000004

$ws.Range("H132").Value = 1451.2858
$ws.Range("I132").Value = 1451.2858
$ws.Range("K132").Value = 4353.857400000001
$ws.Range("M132").Value = -1823.857400000001

$ws.Range("H134").Value = 32714.5
$ws.Range("J134").Value = 32714.5
$ws.Range("L134").Value = 98143.5
$ws.Range("N134").Value = -103213.5

$ws.Range("H136").Value = 2313.6829
$ws.Range("I136").Value = 1617.7878
$ws.Range("J136").Value = 5184.25
$ws.Range("K136").Value = 4853.3634
$ws.Range("L136").Value = 15552.75
$ws.Range("M136").Value = -2303.3634
$ws.Range("N136").Value = -20652.75
